# Add "Swiss" market test-data worksheet (after "Czech"), matching the
# existing Germany/Belgium/Czech sheets' layout.

$wb = $excel.ActiveWorkbook

$czech = $wb.Worksheets.Item("Czech")

# Select the whole Czech sheet before switching away from it, so its
# stored selection becomes the full-sheet range (matches what happens in
# Excel when a user selects everything before moving to a new tab).
$czech.Activate()
$czech.Cells.Select()

# Duplicate the Czech sheet to seed the new one - this keeps every column
# width, merged cell, style and page setting identical, which is exactly
# how the Germany/Belgium/Czech sheets already relate to one another.
$czech.Copy($null, $czech)
$swiss = $wb.Worksheets.Item("Czech (2)")
$swiss.Name = "Swiss"

# Market header + user-story reference.
$swiss.Range("B2").Value = "Switzerland Market"

# Replace the product list (column A, rows 8-34) with the Swiss lineup.
$swiss.Range("A8").Value = "Black Box"
$swiss.Range("A9").Value = "MX1000"
$swiss.Range("A10").Value = "MX4000"
$swiss.Range("A11").Value = "MZX250CH"
$swiss.Range("A12").Value = "MZX251"
$swiss.Range("A13").Value = "MZX252CH"
$swiss.Range("A14").Value = "MZX254CH"
$swiss.Range("A15").Value = "P115D"
$swiss.Range("A16").Value = "P115S"
$swiss.Range("A17").Value = "P115SCH"
$swiss.Range("A18").Value = "P405D"
$swiss.Range("A19").Value = "P405DCH"
$swiss.Range("A20").Value = "P485D"
$swiss.Range("A21").Value = "P805D"
$swiss.Range("A22").Value = "P885D"
$swiss.Range("A23").Value = "Pro16xBB"
$swiss.Range("A24").Value = "Pro16xD"
$swiss.Range("A25").Value = "Pro215D"
$swiss.Range("A26").Value = "Pro215S"
$swiss.Range("A27").Value = "Pro215S Lite CH"
$swiss.Range("A28").Value = "Pro32xBB"
$swiss.Range("A29").Value = "Pro32xD"
$swiss.Range("A30").Value = "Pro415S Lite CH"
$swiss.Range("A31").Value = "Pro815D"
$swiss.Range("A32").Value = "Pro815D-CH"
$swiss.Range("A33").Value = "Pro885D"
$swiss.Range("A34").Value = "ZX1"

# User story reference value (set last, like the original edit).
$swiss.Range("B4").Value = "NGC-3476/T2343/T2641"

# Leave the cursor on B5, same as the source edit.
$swiss.Range("B5").Select()
